$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits (student numbers reshuffled/changed, names partially masked for privacy) ---

# Row 2: B2 student number change, C2 name masked (鍾昊良 -> 鍾x良)
$ws.Range("B2").Value = 108236011
$ws.Range("C2").Value = "鍾x良"

# Row 3: unchanged

# Row 4: B4 student number change
$ws.Range("B4").Value = 132406010

# Row 5: B5 student number change
$ws.Range("B5").Value = 108326016

# Row 6: B6 student number change
$ws.Range("B6").Value = 108546004

# Row 7: B7 student number change, C7 name masked (林育緯 -> 林x緯)
$ws.Range("B7").Value = 107326027
$ws.Range("C7").Value = "林x緯"

# Row 8: B8 student number change, C8 name masked (張峰偉 -> 張x偉)
$ws.Range("B8").Value = 108443005
$ws.Range("C8").Value = "張x偉"

# Row 9: B9 student number change, C9 name masked (王建勳 -> 王x勳)
$ws.Range("B9").Value = 108426020
$ws.Range("C9").Value = "王x勳"

# Row 10: C10 name masked (陳科多 -> 陳x多)
$ws.Range("C10").Value = "陳x多"

# Row 11: B11 student number change, C11 name masked (周賢裕 -> 周x裕)
$ws.Range("B11").Value = 110416001
$ws.Range("C11").Value = "周x裕"

# --- View / formatting edits ---

# Column widths (re-fit to new content)
$ws.Columns("A").ColumnWidth = 9.714285714285715
$ws.Columns("B").ColumnWidth = 10
$ws.Columns("C").ColumnWidth = 7.571428571428572

# Zoom + selection, as left by the editing session
$excel.ActiveWindow.Zoom = 190
$ws.Range("B6").Select() | Out-Null
